# Dodanie podziału treningu na części
# Adds a "Trening" (training) column that splits the session into parts
# (Duza Gra / Mala Gra), converts the date column to a real date value,
# and duplicates rows 2-3 into new rows 4-5 for the second part of the
# training.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header in F1
$ws.Range("F1").Value = "Trening"

# Column A: turn the "27.01.2025" text into a real Excel date serial value.
# The format is first set with a lowercase format code and then corrected
# to the uppercase one (matching the original edit history), and the
# resulting number format is then propagated to the other date cells so
# that every cell shares the very same style.
$ws.Range("A2").Value = 45684
$ws.Range("A2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$dateFormat = $ws.Range("A2").NumberFormat

$ws.Range("A3").Value = 45684
$ws.Range("A3").NumberFormat = $dateFormat

$ws.Range("A4").Value = 45684
$ws.Range("A4").NumberFormat = $dateFormat

$ws.Range("A5").Value = 45684
$ws.Range("A5").NumberFormat = $dateFormat

# New rows 4 and 5 mirror rows 2 and 3 (B/C/D stay empty, E repeats the
# Velocity_Bin pattern) but belong to the "Mala Gra" part of the training.
$ws.Range("E4").Value = "10-15"
$ws.Range("E5").Value = "5-10"

# Tag every data row with which part of the training it belongs to.
$ws.Range("F2").Value = "Duża Gra"
$ws.Range("F3").Value = "Duża Gra"
$ws.Range("F4").Value = "Mała Gra"
$ws.Range("F5").Value = "Mała Gra"
